$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each price/volume/name/link cell below is refreshed to the latest scraped value.
# We assign through a literal text formula (=""value"") and then immediately
# copy/paste-special the cell back onto itself as a value. That final step
# guarantees Excel keeps the data as plain text (preserving formatting such as
# "1.00", "0.0000175" or "63.238.29") instead of silently re-interpreting it as
# a floating point number, while leaving no formula or new cell style behind.
$cell = $ws.Range("D2")
$cell.Formula = "=""63.238.29"""
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("E2")
$cell.Formula = "=""  -1.82%  """
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("D3")
$cell.Formula = "=""3.422.52"""
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("E3")
$cell.Formula = "=""  -2.40%  """
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("D4")
$cell.Formula = "=""1.00"""
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("E4")
$cell.Formula = "=""  -0.01%  """
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("D5")
$cell.Formula = "=""578.32"""
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("E5")
$cell.Formula = "=""  -1.92%  """
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("D6")
$cell.Formula = "=""127.65"""
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("E6")
$cell.Formula = "=""  -5.22%  """
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("E7")
$cell.Formula = "=""  +0.04%  """
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("D8")
$cell.Formula = "=""3.421.33"""
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("E8")
$cell.Formula = "=""  -2.39%  """
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("E9")
$cell.Formula = "=""  -2.03%  """
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("D10")
$cell.Formula = "=""7.50"""
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("E10")
$cell.Formula = "=""  +0.25%  """
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("D11")
$cell.Formula = "=""0.122"""
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("E11")
$cell.Formula = "=""  -1.42%  """
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("E12")
$cell.Formula = "=""  -1.35%  """
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("D13")
$cell.Formula = "=""4.001.94"""
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("E13")
$cell.Formula = "=""  -2.45%  """
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("E14")
$cell.Formula = "=""  -0.78%  """
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("B15")
$cell.Formula = "=""WrappedEther"""
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("C15")
$cell.Formula = "=""https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"""
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("D15")
$cell.Formula = "=""3.420.50"""
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("E15")
$cell.Formula = "=""  -2.49%  """
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("B16")
$cell.Formula = "=""ShibaInu"""
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("C16")
$cell.Formula = "=""https://coinranking.com/coin/xz24e0BjL+shibainu-shib"""
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("D16")
$cell.Formula = "=""0.0000175"""
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("E16")
$cell.Formula = "=""  -3.91%  """
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("D17")
$cell.Formula = "=""63.214.91"""
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("E17")
$cell.Formula = "=""  -1.81%  """
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("D18")
$cell.Formula = "=""24.93"""
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("E18")
$cell.Formula = "=""  -3.40%  """
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("D19")
$cell.Formula = "=""9.70"""
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("E19")
$cell.Formula = "=""  -2.09%  """
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("D20")
$cell.Formula = "=""5.70"""
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("E20")
$cell.Formula = "=""  -0.97%  """
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("D21")
$cell.Formula = "=""13.20"""
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("E21")
$cell.Formula = "=""  -2.88%  """
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("D22")
$cell.Formula = "=""379.76"""
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("E22")
$cell.Formula = "=""  -3.46%  """
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("D23")
$cell.Formula = "=""0.561"""
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("E23")
$cell.Formula = "=""  -2.40%  """
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("D24")
$cell.Formula = "=""3.558.12"""
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("D25")
$cell.Formula = "=""72.92"""
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("E25")
$cell.Formula = "=""  -2.34%  """
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("E26")
$cell.Formula = "=""  -0.06%  """
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("D27")
$cell.Formula = "=""0.0000109"""
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("E27")
$cell.Formula = "=""  -6.95%  """
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("D28")
$cell.Formula = "=""1.00"""
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("E28")
$cell.Formula = "=""  +0.42%  """
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("B29")
$cell.Formula = "=""RenderToken"""
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("C29")
$cell.Formula = "=""https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"""
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("D29")
$cell.Formula = "=""7.02"""
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("E29")
$cell.Formula = "=""  -5.30%  """
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("B30")
$cell.Formula = "=""PancakeSwap"""
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("C30")
$cell.Formula = "=""https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"""
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("D30")
$cell.Formula = "=""2.17"""
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("E30")
$cell.Formula = "=""  -3.93%  """
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("B31")
$cell.Formula = "=""InternetComputer(DFINITY)"""
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("C31")
$cell.Formula = "=""https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"""
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("D31")
$cell.Formula = "=""7.92"""
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("E31")
$cell.Formula = "=""  -4.33%  """
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("B32")
$cell.Formula = "=""Kaspa"""
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("C32")
$cell.Formula = "=""https://coinranking.com/coin/V8GxkwWow+kaspa-kas"""
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("D32")
$cell.Formula = "=""0.154"""
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("E32")
$cell.Formula = "=""  -2.92%  """
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("D33")
$cell.Formula = "=""1.41"""
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("E33")
$cell.Formula = "=""  -4.67%  """
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("D34")
$cell.Formula = "=""3.450.84"""
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("E34")
$cell.Formula = "=""  -2.23%  """
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("E35")
$cell.Formula = "=""  -0.01%  """
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("D36")
$cell.Formula = "=""22.88"""
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("E36")
$cell.Formula = "=""  -2.30%  """
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("D37")
$cell.Formula = "=""5.35"""
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("E37")
$cell.Formula = "=""  -0.48%  """
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("D38")
$cell.Formula = "=""6.76"""
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("E38")
$cell.Formula = "=""  -2.78%  """
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("B39")
$cell.Formula = "=""Monero"""
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("C39")
$cell.Formula = "=""https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"""
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("D39")
$cell.Formula = "=""164.32"""
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("E39")
$cell.Formula = "=""  -1.79%  """
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("B40")
$cell.Formula = "=""ImmutableX"""
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("C40")
$cell.Formula = "=""https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"""
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("D40")
$cell.Formula = "=""1.51"""
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("E40")
$cell.Formula = "=""  -3.46%  """
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("D41")
$cell.Formula = "=""0.0764"""
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("E41")
$cell.Formula = "=""  -3.14%  """
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("E42")
$cell.Formula = "=""  -3.52%  """
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("E43")
$cell.Formula = "=""  -0.03%  """
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("D44")
$cell.Formula = "=""42.01"""
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("E44")
$cell.Formula = "=""  -0.37%  """
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("D45")
$cell.Formula = "=""4.31"""
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("E45")
$cell.Formula = "=""  -3.02%  """
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("D46")
$cell.Formula = "=""1.60"""
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("E46")
$cell.Formula = "=""  -4.67%  """
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("D47")
$cell.Formula = "=""22.84"""
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("E47")
$cell.Formula = "=""  -9.09%  """
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("E48")
$cell.Formula = "=""  -6.73%  """
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("D49")
$cell.Formula = "=""6.71"""
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("E49")
$cell.Formula = "=""  -1.29%  """
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("D50")
$cell.Formula = "=""0.868"""
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("E50")
$cell.Formula = "=""  -3.62%  """
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("D51")
$cell.Formula = "=""2.261.02"""
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$cell = $ws.Range("E51")
$cell.Formula = "=""  -5.32%  """
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = 0
